$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the Table1 by two rows (23 data rows -> 25 data rows = rows 22..26) ---
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.ListRows.Add()
[void]$tbl.ListRows.Add()

# --- Fill content column-by-column (matches authoring order of the source edit) ---

# Column A - Test Case
$ws.Range("A22").Value() = "Board parameters input validation (very small)"
$ws.Range("A23").Value() = "Board parameters input validation (very large)"
$ws.Range("A24").Value() = "Board parameters input validation (floats)"
$ws.Range("A25").Value() = "(User specified) Ant coordinates input validation (floats)"
$ws.Range("A26").Value() = "(User specified) Ant coordinates input validation (on board's border)"

# Column B - Input Values
$ws.Range("B22").Value() = "Input = 1"
$ws.Range("B23").Value() = "Input = 1,000"
$ws.Range("B24").Value() = "0 < Input < 1"
$ws.Range("B25").Value() = "0 < Input < 1"
$ws.Range("B26").Value() = "Input = 0;`r`nInput = Columns/Rows"

# Column C - Driver Functions
$ws.Range("C22").Value() = "Board board = Board(rows, columns);`r`nrows/columns = input"
$ws.Range("C23").Value() = "Board board = Board(rows, columns);`r`nrows/columns = input"
$ws.Range("C24").Value() = "Board board = Board(rows, columns);`r`nrows/columns = input"
$ws.Range("C25").Value() = "Ant ant = Ant(rand_x, rand_y, board);`r`nrand_x, rand_y = input"
$ws.Range("C26").Value() = "Ant ant = Ant(rand_x, rand_y, board);`r`nrand_x, rand_y = input"

# Column D - Expected Outcomes
$ws.Range("D22").Value() = "Board prints and program executes as expected"
$ws.Range("D23").Value() = "Board prints and program executes as expected"
$ws.Range("D24").Value() = "Program crashes"
$ws.Range("D25").Value() = "Program crashes"
$ws.Range("D26").Value() = "Loops back to the question, prompting the user for input"

# --- Row heights ---
$ws.Rows.Item(22).RowHeight() = 68
$ws.Rows.Item(23).RowHeight() = 68
$ws.Rows.Item(24).RowHeight() = 68
$ws.Rows.Item(25).RowHeight() = 68
$ws.Rows.Item(26).RowHeight() = 85

# --- Column widths (Test case / Input Values / Driver Functions columns resized) ---
$ws.Columns.Item(2).ColumnWidth() = 13.33
$ws.Columns.Item(3).ColumnWidth() = 32.67

# --- Selection / view state ---
[void]$ws.Range("D27").Select()
